$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("N16").Value = ""

$ws.Range("H17").Value = 945.3516
$ws.Range("J17").Value = 945.3516
$ws.Range("L17").Value = 2836.0548
$ws.Range("N17").Value = -3172.0548

$ws.Range("H86").Value = 3463.7144
$ws.Range("I86").Value = 3707.6667
$ws.Range("K86").Value = 3707.6667
$ws.Range("M86").Value = -2584.6667

$ws.Range("H89").Value = 3463.7144
$ws.Range("I89").Value = 3707.6667
$ws.Range("K89").Value = 18538.3335
$ws.Range("M89").Value = -12922.3335

$ws.Range("H125").Value = 4362.4116
$ws.Range("I125").Value = 2868.6667
$ws.Range("J125").Value = 5177.1816
$ws.Range("K125").Value = 25818.0003
$ws.Range("L125").Value = 46594.6344
$ws.Range("M125").Value = -23358.0003
$ws.Range("N125").Value = -51514.6344

$ws.Range("H138").Value = 4160.421
$ws.Range("I138").Value = 3522.7407
$ws.Range("J138").Value = 5725.636
$ws.Range("K138").Value = 10568.2221
$ws.Range("L138").Value = 17176.908
$ws.Range("M138").Value = -5428.222099999999
$ws.Range("N138").Value = -27456.908

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1167.6923
$ws.Range("I2").Value = 773.625
$ws.Range("K2").Value = 773.625
$ws.Range("M2").Value = -660.625

$ws.Range("H32").Value = 1661.8817
$ws.Range("I32").Value = 1359.159
$ws.Range("J32").Value = 6989.8
$ws.Range("K32").Value = 1359.159
$ws.Range("L32").Value = 6989.8
$ws.Range("M32").Value = -1072.159
$ws.Range("N32").Value = -7563.8

$ws.Range("H110").Value = 1750
$ws.Range("I110").Value = 1750
$ws.Range("K110").Value = 1750
$ws.Range("M110").Value = 295

$ws.Range("H116").Value = 1167.6923
$ws.Range("I116").Value = 773.625
$ws.Range("K116").Value = 773.625
$ws.Range("M116").Value = 1520.375

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1167.6923
$ws.Range("I3").Value = 773.625
$ws.Range("K3").Value = 773.625
$ws.Range("M3").Value = -659.625

$ws.Range("H94").Value = 26271
$ws.Range("I94").Value = 1694.6666
$ws.Range("K94").Value = 1694.6666
$ws.Range("M94").Value = -1243.6666

$ws.Range("H99").Value = 2533
$ws.Range("I99").Value = 2749.8333
$ws.Range("K99").Value = 2749.8333
$ws.Range("M99").Value = -1251.8333

$ws.Range("H139").Value = 143101
$ws.Range("J139").Value = 158499.67
$ws.Range("L139").Value = 158499.67
$ws.Range("N139").Value = -168779.67

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 2730216.8
$ws.Range("J19").Value = 1500
$ws.Range("L19").Value = 1500
$ws.Range("N19").Value = -1840

$ws.Range("H24").Value = 2730216.8
$ws.Range("J24").Value = 1500
$ws.Range("L24").Value = 1500
$ws.Range("N24").Value = -1840

$ws.Range("H31").Value = 2967.7693
$ws.Range("I31").Value = 2091.5
$ws.Range("J31").Value = 4055.5518
$ws.Range("K31").Value = 2091.5
$ws.Range("L31").Value = 4055.5518
$ws.Range("M31").Value = -1796.5
$ws.Range("N31").Value = -4645.5518

$ws.Range("H34").Value = 2967.7693
$ws.Range("I34").Value = 2091.5
$ws.Range("J34").Value = 4055.5518
$ws.Range("K34").Value = 2091.5
$ws.Range("L34").Value = 4055.5518
$ws.Range("M34").Value = -1889.5
$ws.Range("N34").Value = -4459.5518

$ws.Range("H74").Value = 42500
$ws.Range("J74").Value = 42500
$ws.Range("L74").Value = 42500
$ws.Range("N74").Value = -44248

$ws.Range("H77").Value = 42500
$ws.Range("J77").Value = 42500
$ws.Range("L77").Value = 127500
$ws.Range("N77").Value = -136236

$ws.Range("H107").Value = 976.73334
$ws.Range("I107").Value = 329.76923
$ws.Range("K107").Value = 329.76923
$ws.Range("M107").Value = 1590.23077

$ws.Range("H132").Value = 3715.1292
$ws.Range("I132").Value = 3535.45
$ws.Range("K132").Value = 10606.35
$ws.Range("M132").Value = -8076.349999999999

$ws.Range("H134").Value = 3321.5898
$ws.Range("I134").Value = 3191.8276
$ws.Range("J134").Value = 3697.9
$ws.Range("K134").Value = 9575.4828
$ws.Range("L134").Value = 11093.7
$ws.Range("M134").Value = -7040.4828
$ws.Range("N134").Value = -16163.7

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 709069.0600000001
$ws.Range("I11").Value = 784944.3
$ws.Range("K11").Value = 2354832.9
$ws.Range("M11").Value = -2354692.9

$ws.Range("H39").Value = 4229.3335
$ws.Range("I39").Value = 2466.3333
$ws.Range("J39").Value = 4523.1665
$ws.Range("K39").Value = 7398.999899999999
$ws.Range("L39").Value = 13569.4995
$ws.Range("M39").Value = -7104.999899999999
$ws.Range("N39").Value = -14157.4995

$ws.Range("H104").Value = 4833.1665
$ws.Range("I104").Value = 4199.8
$ws.Range("J104").Value = 8000
$ws.Range("K104").Value = 12599.4
$ws.Range("L104").Value = 24000
$ws.Range("M104").Value = -9978.400000000001
$ws.Range("N104").Value = -29242

$ws.Range("H109").Value = 2251.0667
$ws.Range("I109").Value = 1697.5714
$ws.Range("J109").Value = 10000
$ws.Range("K109").Value = 5092.7142
$ws.Range("L109").Value = 30000
$ws.Range("M109").Value = -4052.7142
$ws.Range("N109").Value = -32080

$ws.Range("H113").Value = 1795.8889
$ws.Range("I113").Value = 1489.25
$ws.Range("J113").Value = 1883.5
$ws.Range("K113").Value = 4467.75
$ws.Range("L113").Value = 5650.5
$ws.Range("M113").Value = -2297.75
$ws.Range("N113").Value = -9990.5

$ws.Range("H122").Value = 2017
$ws.Range("I122").Value = 432.66666
$ws.Range("J122").Value = 3601.3333
$ws.Range("K122").Value = 3893.99994
$ws.Range("L122").Value = 32411.9997
$ws.Range("M122").Value = -1443.99994
$ws.Range("N122").Value = -37311.9997

$ws.Range("H126").Value = 15510.643
$ws.Range("I126").Value = 7713.375
$ws.Range("J126").Value = 25907
$ws.Range("K126").Value = 23140.125
$ws.Range("L126").Value = 77721
$ws.Range("M126").Value = -18200.125
$ws.Range("N126").Value = -87601

$ws.Range("H131").Value = 7288904.5
$ws.Range("I131").Value = 7408768.5
$ws.Range("J131").Value = 7249818.5
$ws.Range("K131").Value = 22226305.5
$ws.Range("L131").Value = 21749455.5
$ws.Range("M131").Value = -22221265.5
$ws.Range("N131").Value = -21759535.5

$ws.Range("H134").Value = 11419.277
$ws.Range("I134").Value = 5191.0586
$ws.Range("K134").Value = 15573.1758
$ws.Range("M134").Value = -10503.1758

$ws.Range("H141").Value = 22665.523
$ws.Range("I141").Value = 4020.5454
$ws.Range("K141").Value = 12061.6362
$ws.Range("M141").Value = -6881.636200000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 144
$ws.Range("I2").Value = 139
$ws.Range("J2").Value = 159
$ws.Range("K2").Value = 139
$ws.Range("L2").Value = 159
$ws.Range("M2").Value = -26
$ws.Range("N2").Value = -385

$ws.Range("H80").Value = 2841.8572
$ws.Range("I80").Value = 2998
$ws.Range("J80").Value = 2779.4
$ws.Range("K80").Value = 2998
$ws.Range("L80").Value = 2779.4
$ws.Range("M80").Value = -2000
$ws.Range("N80").Value = -4775.4

$ws.Range("H83").Value = 2841.8572
$ws.Range("I83").Value = 2998
$ws.Range("J83").Value = 2779.4
$ws.Range("K83").Value = 14990
$ws.Range("L83").Value = 13897
$ws.Range("M83").Value = -9998
$ws.Range("N83").Value = -23881

$ws.Range("H102").Value = 5379.5
$ws.Range("I102").Value = 5442.857
$ws.Range("J102").Value = 5231.6665
$ws.Range("K102").Value = 5442.857
$ws.Range("L102").Value = 5231.6665
$ws.Range("M102").Value = -3820.857
$ws.Range("N102").Value = -8475.666499999999

$ws.Range("H107").Value = 898.7273
$ws.Range("I107").Value = 801.7143
$ws.Range("K107").Value = 801.7143
$ws.Range("M107").Value = 1118.2857

$ws.Range("H133").Value = 79989.5
$ws.Range("J133").Value = 79989.5
$ws.Range("L133").Value = 79989.5
$ws.Range("N133").Value = -90109.5

$ws.Range("H135").Value = 94084.75
$ws.Range("J135").Value = 94084.75
$ws.Range("L135").Value = 94084.75
$ws.Range("N135").Value = -104224.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3491.3
$ws.Range("I7").Value = 3491.3
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 3491.3
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -3379.3
$ws.Range("N7").Value = ""

$ws.Range("H40").Value = 2242.7222
$ws.Range("I40").Value = 1799.2858
$ws.Range("J40").Value = 2524.9092
$ws.Range("K40").Value = 1799.2858
$ws.Range("L40").Value = 2524.9092
$ws.Range("M40").Value = -1663.2858
$ws.Range("N40").Value = -2796.9092

$ws.Range("H122").Value = 2500
$ws.Range("I122").Value = 2500
$ws.Range("K122").Value = 7500
$ws.Range("M122").Value = -5050

$ws.Range("H126").Value = 3491.3
$ws.Range("I126").Value = 3491.3
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 10473.9
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -8003.900000000001
$ws.Range("N126").Value = ""

$ws.Range("H132").Value = 1981.5883
$ws.Range("I132").Value = 1621.1
$ws.Range("K132").Value = 4863.299999999999
$ws.Range("M132").Value = -2333.299999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 32916.75
$ws.Range("I3").Value = 33333
$ws.Range("J3").Value = 32500.5
$ws.Range("K3").Value = 33333
$ws.Range("L3").Value = 32500.5
$ws.Range("M3").Value = -33219
$ws.Range("N3").Value = -32728.5

$ws.Range("H6").Value = 0
$ws.Range("J6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("N6").Value = ""

$ws.Range("H122").Value = 4777
$ws.Range("I122").Value = 4777
$ws.Range("K122").Value = 14331
$ws.Range("M122").Value = -11881

$ws.Range("H132").Value = 3808
$ws.Range("I132").Value = 4374
$ws.Range("J132").Value = 3242
$ws.Range("K132").Value = 13122
$ws.Range("L132").Value = 9726
$ws.Range("M132").Value = -10592
$ws.Range("N132").Value = -14786
